$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F (shifts F:P -> G:Q) and set its header
$ws.Columns("F:F").Insert()
$ws.Range("F1").Value = "jenis_kepemilikan"

# Rename two headers that moved to J and L after the insert
$ws.Range("J1").Value = "nominal_pembayaran_online"
$ws.Range("L1").Value = "nominal_pembayaran_offline"

# Autofit only the columns whose contents actually changed, so the
# widths of untouched columns are left exactly as they were
$ws.Columns("F:F").AutoFit() | Out-Null
$ws.Columns("J:J").AutoFit() | Out-Null
$ws.Columns("L:L").AutoFit() | Out-Null

# Scroll back to A1 and select H23, matching the saved view state
$ws.Range("A1").Select() | Out-Null
$ws.Range("H23").Select() | Out-Null
